$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Finish out row 7 with the two trailing columns that were missing.
$ws.Range("X7").Value = -3.5699769999999944
$ws.Range("Y7").Value = "Down"

# Append a brand-new row 8 of data (mirrors the structure of the prior rows).
$ws.Range("A8").Value = 42649.8909375
$ws.Range("A8").NumberFormat = "m/d/yyyy h:mm"

$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "Neutral"
$ws.Range("D8").Value = 36
$ws.Range("E8").Value = 21742
$ws.Range("F8").Value = 1319
$ws.Range("G8").Value = 56
$ws.Range("H8").Value = 40
$ws.Range("I8").Value = 89
$ws.Range("J8").Value = 9
$ws.Range("K8").Value = 37852
$ws.Range("L8").Value = 267
$ws.Range("M8").Value = 193
$ws.Range("N8").Value = 131
$ws.Range("O8").Value = 14
$ws.Range("P8").Value = "Named"
$ws.Range("Q8").Value = 49.72799223503381
$ws.Range("R8").Value = 0

$ws.Range("S8").Value = 0.1095
$ws.Range("S8").NumberFormat = "0.00%"
$ws.Range("T8").Value = 0.0025000000000000001
$ws.Range("T8").NumberFormat = "0.00%"

$ws.Range("U8").Value = 5.95
$ws.Range("V8").Value = "N/A"
$ws.Range("W8").Value = 0
